$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 77.63718533333333
$ws.Range("H2").Value = 232.911556
$ws.Range("I2").Value = 0.4513549673384918
$ws.Range("J2").Value = 0.4513549673384918
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.465868666666666
$ws.Range("N2").Value = 13.397606
$ws.Range("O2").Value = 0.3171934613136183
$ws.Range("P2").Value = 0.3171934613136182
$ws.Range("Q2").Value = 346.7174733483262
$ws.Range("R2").Value = 3120.457260134936
$ws.Range("S2").Value = 0.1431668443711913
$ws.Range("T2").Value = 0.1431668443711913

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 77.63718533333333
$ws.Range("H3").Value = 232.911556
$ws.Range("I3").Value = 0.4513549673384918
$ws.Range("J3").Value = 0.4513549673384918
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.701621333333333
$ws.Range("N3").Value = 20.104864
$ws.Range("O3").Value = 0.4759903673387287
$ws.Range("P3").Value = 0.4759903673387287
$ws.Range("Q3").Value = 520.2950174898205
$ws.Range("R3").Value = 4682.655157408383
$ws.Range("S3").Value = 0.2148406167036087
$ws.Range("T3").Value = 0.2148406167036086

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 77.63718533333333
$ws.Range("H4").Value = 232.911556
$ws.Range("I4").Value = 0.4513549673384918
$ws.Range("J4").Value = 0.4513549673384918
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.911831333333334
$ws.Range("N4").Value = 8.735494000000001
$ws.Range("O4").Value = 0.2068161713476531
$ws.Range("P4").Value = 0.206816171347653
$ws.Range("Q4").Value = 226.0663888854072
$ws.Range("R4").Value = 2034.597499968664
$ws.Range("S4").Value = 0.09334750626369188
$ws.Range("T4").Value = 0.09334750626369184

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 62.40815866666667
$ws.Range("H5").Value = 187.224476
$ws.Range("I5").Value = 0.3628188257432201
$ws.Range("J5").Value = 0.3628188257432201
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.465868666666666
$ws.Range("N5").Value = 13.397606
$ws.Range("O5").Value = 0.3171934613136183
$ws.Range("P5").Value = 0.3171934613136182
$ws.Range("Q5").Value = 278.7066403338284
$ws.Range("R5").Value = 2508.359763004456
$ws.Range("S5").Value = 0.1150837591672345
$ws.Range("T5").Value = 0.1150837591672345

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 62.40815866666667
$ws.Range("H6").Value = 187.224476
$ws.Range("I6").Value = 0.3628188257432201
$ws.Range("J6").Value = 0.3628188257432201
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.701621333333333
$ws.Range("N6").Value = 20.104864
$ws.Range("O6").Value = 0.4759903673387287
$ws.Range("P6").Value = 0.4759903673387287
$ws.Range("Q6").Value = 418.2358474945849
$ws.Range("R6").Value = 3764.122627451264
$ws.Range("S6").Value = 0.1726982661429216
$ws.Range("T6").Value = 0.1726982661429215

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 62.40815866666667
$ws.Range("H7").Value = 187.224476
$ws.Range("I7").Value = 0.3628188257432201
$ws.Range("J7").Value = 0.3628188257432201
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.911831333333334
$ws.Range("N7").Value = 8.735494000000001
$ws.Range("O7").Value = 0.2068161713476531
$ws.Range("P7").Value = 0.206816171347653
$ws.Range("Q7").Value = 181.7220318612383
$ws.Range("R7").Value = 1635.498286751144
$ws.Range("S7").Value = 0.07503680043306409
$ws.Range("T7").Value = 0.07503680043306407

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 31.96380833333333
$ws.Range("H8").Value = 95.891425
$ws.Range("I8").Value = 0.1858262069182881
$ws.Range("J8").Value = 0.1858262069182881
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.465868666666666
$ws.Range("N8").Value = 13.397606
$ws.Range("O8").Value = 0.3171934613136183
$ws.Range("P8").Value = 0.3171934613136182
$ws.Range("Q8").Value = 142.7461701031722
$ws.Range("R8").Value = 1284.71553092855
$ws.Range("S8").Value = 0.05894285777519243
$ws.Range("T8").Value = 0.05894285777519241

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 31.96380833333333
$ws.Range("H9").Value = 95.891425
$ws.Range("I9").Value = 0.1858262069182881
$ws.Range("J9").Value = 0.1858262069182881
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.701621333333333
$ws.Range("N9").Value = 20.104864
$ws.Range("O9").Value = 0.4759903673387287
$ws.Range("P9").Value = 0.4759903673387287
$ws.Range("Q9").Value = 214.2093398212444
$ws.Range("R9").Value = 1927.8840583912
$ws.Range("S9").Value = 0.08845148449219856
$ws.Range("T9").Value = 0.08845148449219853

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 31.96380833333333
$ws.Range("H10").Value = 95.891425
$ws.Range("I10").Value = 0.1858262069182881
$ws.Range("J10").Value = 0.1858262069182881
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.911831333333334
$ws.Range("N10").Value = 8.735494000000001
$ws.Range("O10").Value = 0.2068161713476531
$ws.Range("P10").Value = 0.206816171347653
$ws.Range("Q10").Value = 93.07321863766111
$ws.Range("R10").Value = 837.6589677389501
$ws.Range("S10").Value = 0.0384318646508971
$ws.Range("T10").Value = 0.03843186465089708
